$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.708.36'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '3.027.54'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '583.03'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.62'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.56%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -2.25%  '
$ws.Range('D9').Value = '3.027.94'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('E10').Value = '  -2.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.70'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.60%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.444'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.24%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000231'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.22%  '
$ws.Range('E14').Value = '  -3.94%  '
$ws.Range('E15').Value = '  +2.29%  '
$ws.Range('D16').Value = '3.526.55'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.07'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = '62.678.83'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = '3.024.89'
$ws.Range('E19').Value = '  -1.39%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '468.65'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.84%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.05'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.88%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.692'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.48%  '
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '80.97'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.41'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.57%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.35'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.26%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('B29').Value = 'NEARProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.28'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.53%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('E31').Value = '  -0.63%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.14'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.46'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.74%  '
$ws.Range('E34').Value = '  -4.05%  '
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('D36').Value = '0.0₃0796'
$ws.Range('E36').Value = '  -2.96%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.78'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.86%  '
$ws.Range('E38').Value = '  -1.52%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '50.24'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '9.02'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.92%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.94'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -12.12%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '424.37'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.77%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.282'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.114'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('D45').Value = '2.802.13'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0356'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '37.87'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -7.61%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '129.08'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '24.37'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.44%  '
$ws.Range('E51').Value = '  -0.43%  '
